# Refresh the cryptos list (price + 1h volume change columns, and the
# Stellar/Stacks row ordering swap) to match the latest GitHub Actions run.
#
# Note: several "Price" values look like plain numbers (e.g. "213.98",
# "0.0000220"). Assigning those to Range.Value directly would make Excel's
# COM layer auto-convert them to real numeric cells (losing formatting like
# trailing zeros, e.g. "0.0000220" -> 2.2E-05). To keep them as text - as
# they are in the source workbook - we prefix those with a leading
# apostrophe, which is the standard Excel trick to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.156.23"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "3.117.88"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'213.98"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").Value = "'635.79"
$ws.Range("E6").Value = "  +3.89%  "

$ws.Range("D7").Value = "'0.386"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").Value = "'0.823"
$ws.Range("E8").Value = "  +20.19%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "3.112.30"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "'5.43"
$ws.Range("E14").Value = "  +4.13%  "

$ws.Range("D15").Value = "88.008.45"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").Value = "3.690.72"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "'32.03"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "3.100.82"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("E19").Value = "  +3.06%  "

$ws.Range("D20").Value = "'0.0000220"
$ws.Range("E20").Value = "  +18.18%  "

$ws.Range("D21").Value = "'13.24"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").Value = "'423.96"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").Value = "'8.47"
$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("E24").Value = "  -3.08%  "

$ws.Range("D25").Value = "'5.55"
$ws.Range("E25").Value = "  +9.28%  "

$ws.Range("D26").Value = "'83.01"
$ws.Range("E26").Value = "  +10.99%  "

$ws.Range("D27").Value = "'11.51"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").Value = "3.276.33"
$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -6.15%  "

$ws.Range("D32").Value = "'3.92"
$ws.Range("E32").Value = "  -3.95%  "

$ws.Range("E33").Value = "  -2.67%  "

$ws.Range("D34").Value = "'0.149"
$ws.Range("E34").Value = "  +17.64%  "

$ws.Range("D35").Value = "'500.74"
$ws.Range("E35").Value = "  -5.06%  "

$ws.Range("D36").Value = "'6.78"
$ws.Range("E36").Value = "  -2.67%  "

$ws.Range("E37").Value = "  +1.40%  "

$ws.Range("D38").Value = "'1.84"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").Value = "'22.38"
$ws.Range("E39").Value = "  +2.86%  "

$ws.Range("D40").Value = "'22.21"
$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "'0.367"
$ws.Range("E43").Value = "  -0.73%  "

# Row 44/45: Stellar and Stacks swapped places (ranking order changed),
# while the index column (A) stays as-is.
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.84"
$ws.Range("E44").Value = "  -3.18%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.137"
$ws.Range("E45").Value = "  +11.30%  "

$ws.Range("D46").Value = "'146.05"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("D47").Value = "'43.91"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "'0.0656"
$ws.Range("E48").Value = "  +12.24%  "

$ws.Range("D49").Value = "'162.74"
$ws.Range("E49").Value = "  -4.81%  "

$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("E51").Value = "  -2.76%  "
